$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# Rename empires (B column)
$ws.Cells.Item(35, 2).Value = "Королевство Леон"
$ws.Cells.Item(36, 2).Value = "Королевство Англия"
$ws.Cells.Item(38, 2).Value = "Королевство Дания"
$ws.Cells.Item(40, 2).Value = "Империя Сун"
$ws.Cells.Item(41, 2).Value = "Королевство Швеция"
$ws.Cells.Item(43, 2).Value = "Королевство Венгрия"
$ws.Cells.Item(44, 2).Value = "Королевство Франция"

# Update population values (D21, D22)
Set-TextValue 21 4 "0.065"
Set-TextValue 22 4 "0.080"

# Add description/population/capital for rows 34-44
$ws.Cells.Item(34, 3).Value = "Государство на Корейском полуострове, появившееся после падения государства Силла в 935 году и существовавшее до воцарения династии Чосон в 1392 году."
Set-TextValue 34 4 "4.940"
$ws.Cells.Item(34, 5).Value = "Кэгён"

$ws.Cells.Item(35, 3).Value = "Средневековое королевство в северо-западной части Пиренейского полуострова."
Set-TextValue 35 4 "2.220"
$ws.Cells.Item(35, 5).Value = "Леон"

$ws.Cells.Item(36, 3).Value = "Государство в Западной Европе, занимавшее южную часть острова Великобритания, известное с 927 года. После объединения с королевством Шотландия в 1707 году трансформировалось в королевство Великобритания."
Set-TextValue 36 4 "4.345"
$ws.Cells.Item(36, 5).Value = "Винчестер"

$ws.Cells.Item(37, 3).Value = "Возникло после объединения в 933 году королевств Верхней и Нижней Бургундии. В 1032—1034 годах государство после смерти бездетного короля Рудольфа III и последовавшей за ней войны за Бургундское наследство вошло в состав Священной Римской империи в качестве третьего королевства, наряду с Германией и Италией."
Set-TextValue 37 4 "1.780"
$ws.Cells.Item(37, 5).Value = "Арль"

$ws.Cells.Item(38, 3).Value = "Государство монархическое. Зародилось на территории современной Дании в результате слияния отдельных племенных групп в одно территориальное государство."
Set-TextValue 38 4 "1.270"
$ws.Cells.Item(38, 5).Value = "Копенгаген"

$ws.Cells.Item(39, 3).Value = "В 877 году после завоевания Малой Польши Великой Моравией центром формирования польского государства осталась Великая Польша."
Set-TextValue 39 4 "3.650"
$ws.Cells.Item(39, 5).Value = "Гнезно"

$ws.Cells.Item(40, 3).Value = "Государство в Китае, существовавшее с 960 по 1279 год. Правящая династия — Чжао, по фамилии рода государей. Основание империи положило конец раздробленности Китая, продолжавшейся со времени падения империи Тан в 907 году. "
Set-TextValue 40 4 "100.000"
$ws.Cells.Item(40, 5).Value = " Бяньлян"

$ws.Cells.Item(41, 3).Value = "Государство монархическое. Зародилось на территории современной Швеции в результате объединения племен свеев и гетов родом Инглингов."
Set-TextValue 41 4 "1.940"
$ws.Cells.Item(41, 5).Value = "Упсала"

$ws.Cells.Item(42, 3).Value = "Надгосударственный союз немецких, итальянских, балканских, франкских и западнославянских государств и народов, существовавший на протяжении 844 лет."
Set-TextValue 42 4 "18.860"
$ws.Cells.Item(42, 5).Value = "Ахен"

$ws.Cells.Item(43, 3).Value = "Зародилось на терр. совр. Венгрии в результате нашествия венгерских племен и образования Королевства Венгрия."
Set-TextValue 43 4 "4.680"
$ws.Cells.Item(43, 5).Value = "Эстергом"

$ws.Cells.Item(44, 3).Value = "Государство (королевство) в Западной Европе, предшественник современной Французской республики.`nФранцузское королевство одно из самых мощных государств в Западной Европе, великая держава со времён позднего средневековья и Столетней войны. Была также ранней колониальной державой, со значительными владениями в Северной Америке."
$ws.Cells.Item(44, 5).Value = "Париж"

# Update sheet view
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("B45").Select()